$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-05-24 Saturday" "2025-05-25 Sunday"

Replace-Text "744×5=" "420×7="
Replace-Text "986×7=" "965×5="
Replace-Text "876×4=" "330×4="
Replace-Text "899×4=" "459×6="
Replace-Text "570×8=" "258×3="

Replace-Text "382×7=" "335×5="
Replace-Text "221×4=" "470×7="
Replace-Text "503×5=" "462×7="
Replace-Text "889×8=" "354×9="
Replace-Text "709×6=" "336×3="

Replace-Text "269×7=" "656×9="
Replace-Text "138×3=" "511×4="
Replace-Text "529×3=" "495×5="
Replace-Text "663×8=" "139×2="
Replace-Text "277×9=" "411×9="

Replace-Text "367×2=" "744×3="
Replace-Text "950×5=" "888×8="
Replace-Text "451×3=" "263×8="
Replace-Text "428×3=" "447×8="
Replace-Text "399×9=" "390×7="

Replace-Text "897×3=" "177×7="
Replace-Text "386×8=" "860×4="
Replace-Text "514×4=" "909×7="
Replace-Text "198×3=" "656×6="
Replace-Text "209×9=" "332×9="
